$d = $word.ActiveDocument

# 1. Add " №2" (as two separate runs) to the end of the Title paragraph,
#    right after the existing "работе" run, so that the title reads
#    "Отчёт по лабораторной работе №2".
$titlePara = $d.Paragraphs(1)

$afterWork = $titlePara.Range
$afterWork.InsertAfter(" ")

$afterSpace = $titlePara.Range
$afterSpace.InsertAfter("№2")

# 2. Fix the typo in the author's patronymic: "Риантвоич" -> "Ринатович".
$d.Content.Find.Execute("Риантвоич", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Ринатович", 2)
